$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, pushing the existing rows 164..178 down to 165..179
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record
$ws.Cells.Item(164, 1).Value = 4
$ws.Cells.Item(164, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(164, 3).Value = "Los Lagos"
$ws.Cells.Item(164, 4).Value = 44449
$ws.Cells.Item(164, 5).Value = 10
$ws.Cells.Item(164, 6).Value = 100114014
$ws.Cells.Item(164, 7).Value = "Betarraga"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 1000
$ws.Cells.Item(164, 11).Value = 1000
$ws.Cells.Item(164, 12).Value = 1000
$ws.Cells.Item(164, 13).Value = 1000
$ws.Cells.Item(164, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 200
$ws.Cells.Item(164, 17).Value = 5
$ws.Cells.Item(164, 18).Value = "Hortaliza"
